$wb = $excel.ActiveWorkbook

# --- Servicios de Salud sheet: rework the "institutions of health affiliation" table header ---
$ws = $wb.Worksheets.Item("Servicios de Salud")

# Unmerge the A5:A7 block that used to hold the "Intituciones de afiliaciones..." label
$ws.Range("A5:A7").UnMerge()

# Move the label text down from A5 to A7 (set the new location BEFORE clearing the old
# one so the shared string stays referenced), then clear the old cell entirely
$label = $ws.Range("A5").Value2
$ws.Range("A7").Value2 = $label
$ws.Range("A5").Clear()

# Re-style A6 and A7 with a new format: no border, vertically centered, wrapped text
$ws.Range("A6:A7").Borders.LineStyle = -4142
$ws.Range("A6:A7").WrapText = $true
$ws.Range("A6:A7").VerticalAlignment = -4108
$ws.Range("A6:A7").HorizontalAlignment = 1

# Adjust row heights: row 5 keeps the default (12pt) explicitly, row 7 grows for wrapped text
$ws.Rows.Item(5).RowHeight = 12
$ws.Rows.Item(7).RowHeight = 24

# Replace the generic Total/Hombres/Mujeres header labels (row 7) with distinct ones
$ws.Range("B7").Value = "Total_SD"
$ws.Range("C7").Value = "Hombres_SD"
$ws.Range("D7").Value = "Mujeres_SD"
$ws.Range("E7").Value = "Total_CD"
$ws.Range("F7").Value = "Hombres_CD"
$ws.Range("G7").Value = "Mujeres_CD"

# Select A7 and make "Servicios de Salud" the active/selected tab
$ws.Range("A7").Select()
$ws.Activate()
